$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column A updates: some values become text strings, others remain/become numbers
$ws.Range("A2").Value = "qqqq1"
$ws.Range("A3").Value = 343442
$ws.Range("A4").Value = "yuyuyu666"
$ws.Range("A7").Value = "686868efjf"
$ws.Range("A8").Value = "fkrioi5"
$ws.Range("A9").Value = "oti5o4oi"
$ws.Range("A10").Value = "oi5pir"
$ws.Range("A11").Value = 767
